# Add three new character styles (GaNStyle, GaNParagraph, GaNLinks) and
# apply them to the runs of text that were newly styled in the commit.

$d = $word.ActiveDocument

# wdStyleTypeCharacter = 2
$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608   # RGB(0, 0, 128) -> stored as w:color="000080"
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1     # wdUnderlineSingle

# 1) The "Perioadele campaniei ..." run appears four times in the document;
#    style every occurrence with GaNStyle.
$searchText = "Perioadele campaniei din 2022 pentru Constelația Bootes: 14-23 mai, 13-22 iunie, 12-21 iulie"
$rng = $d.Content
while ($rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $rng.End = $d.Content.End
}

# 2) The "Prin această activitate ..." paragraph run gets GaNParagraph.
$paraText = "Prin această activitate participați în cadrul unei campanii globale de observare și consemnare a celor mai slabe stele vizibile ca metodă de măsurare a poluării luminoase dintr-un anumit loc. Localizând și observând  Constelația Bootes pe cerul nopții și comparând-o cu diagramele stelare, oamenii din întreaga lume vor putea afla în ce măsură iluminatul nocturn din comunitatea lor contribuie la poluarea luminoasă. Contribuțiile dumneavoastră la baza de date online vor facilita o documentare globală privind cerul nocturn observabil."
$rng2 = $d.Content
$rng2.Find.Execute($paraText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Style = "GaNParagraph"

# 3) The "de Jan Hollan, CzechGlobe ..." run gets GaNLinks.
$linkText = "de Jan Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3 = $d.Content
$rng3.Find.Execute($linkText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3.Style = "GaNLinks"

Write-Output "Styles added and applied."
